$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the DNP / mechanical placeholder rows that were dropped from the
# purchasing BOM (H1,H2,H4 / H3 / H5,H6 / J3 / R24 / TP.../TP9). Deleting
# from the bottom up so earlier row numbers stay valid while we work.
$rowsToDelete = @(32, 31, 24, 19, 16, 15, 14)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# The rows that used to be U1..Y1 (now rows 26-31 after the shift) didn't
# carry a "Reference" value in column J before; the refreshed sheet adds it.
for ($r = 26; $r -le 31; $r++) {
    $ws.Cells.Item($r, 10).Value = "PGKB"
}

# Update the view state left after the edit (active cell / scroll).
$ws.Range("B20").Select()
